$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (Michel Foucault) previously only had the philosopher name/initial
# (columns A and B). This adds the set of quotes across columns C..P,
# matching the shared-string table growth (uniqueCount 270 -> 284).
#
# A handful of cells (D12, E12, H12, J12) pick up the same
# left/center/indent-1 cell style ("s=\"1\"") that many other quote cells
# in the sheet already use; we grab that formatting from an existing cell
# (C2) via copy/paste-special so we reuse the existing style index instead
# of fabricating a new one.
$ws.Range("C2").Copy() | Out-Null

$ws.Range("C12").Value = """Oamenii știu ce fac; adesea știu de ce fac ceea ce fac; dar ceea ce nu știu este ce efect are ceea ce fac"""
$ws.Range("D12").Value = """Nu simt că este necesar să știu exact ce sunt. Interesul principal în viață și muncă este să devii altcineva față de ceea ce erai la început"""
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = """Unde există putere, există rezistență"""
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("F12").Value = """De ce ar trebui lampa sau casa să fie un obiect de artă, dar nu și viața noastră?"""
$ws.Range("G12").Value = """Nu sunt un profet. Meseria mea este să fac ferestre acolo unde erau odată pereți"""
$ws.Range("H12").Value = """Cunoașterea nu este pentru a ști: cunoașterea este pentru a tăia"""
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Value = """Poate că ținta în zilele noastre nu este să descoperim ce suntem, ci să refuzăm ceea ce suntem"""
$ws.Range("J12").Value = """Școlile îndeplinesc aceleași funcții sociale ca și închisorile și instituțiile de sănătate mintală - de a defini, clasifica, controla și reglementa oamenii"""
$ws.Range("J12").PasteSpecial(-4122) | Out-Null
$ws.Range("K12").Value = """Nu mă întreba cine sunt și nu mă ruga să rămân la fel. Mai mult de o persoană, fără îndoială asemănătoare mie, scrie pentru a nu avea un chip"""
$ws.Range("L12").Value = """Iluminismul, este cel care a descoperit libertățile, a inventat și disciplinle"""
$ws.Range("M12").Value = """Ce dorință poate fi contrară naturii, având în vedere că a fost dăruită omului de către natură însăși?"""
$ws.Range("N12").Value = """Nu mă întreba cine sunt și nu mă ruga să rămân la fel"""
$ws.Range("O12").Value = """Nu credeți că trebuie să fii trist pentru a fi militant, chiar dacă ceea ce combați este abominabil"""
$ws.Range("P12").Value = """Moartea a părăsit vechiul ei cer tragic și a devenit nucleul liric al omului: adevărul său invizibil, secretul său vizibil"""

# Formatting has been pasted; drop out of copy mode.
$excel.CutCopyMode = 0

# Match the saved selection/active cell from the edited workbook.
$ws.Range("Q12").Select() | Out-Null
